$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

$ws1.Range("A2").Value = "96555b91-741f-48b4-9887-5c2f343ac0d9.md"
$ws1.Range("B2").Value = "e2e\96555b91-741f-48b4-9887-5c2f343ac0d9.md"
$ws1.Range("G2").Value = "2016-08-25 09:02:11"

$ws2.Range("A2").Value = "96555b91-741f-48b4-9887-5c2f343ac0d9.md"
$ws2.Range("G2").Value = "96555b91-741f-48b4-9887-5c2f343ac0d9.c39b27a6ee30c08d8156d4f335606b59bf26ae0d.zh-cn.xlf"
$ws2.Range("H2").Value = "2016-08-25 09:01:58"

$ws3.Range("A2").Value = "96555b91-741f-48b4-9887-5c2f343ac0d9.md"
$ws3.Range("G2").Value = "96555b91-741f-48b4-9887-5c2f343ac0d9.c39b27a6ee30c08d8156d4f335606b59bf26ae0d.de-de.xlf"
$ws3.Range("H2").Value = "2016-08-25 09:02:11"
